$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell holding the default (unstyled) format, used to strip the
# NumberFormat-driven style churn that Range.Value induces when a literal
# numeric-looking string (e.g. "1.009") is assigned to a General cell.
$styleDonor = $ws.Range("D6")

$ws.Range("D2").Value = '29.489.28'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '1.913.86'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = $styleDonor.Style
$ws.Range("E4").Value = '  +0.63%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.77'
$ws.Range("D5").Style = $styleDonor.Style
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4823'
$ws.Range("D7").Style = $styleDonor.Style
$ws.Range("E7").Value = '  +1.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4070'
$ws.Range("D8").Style = $styleDonor.Style
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08147'
$ws.Range("D9").Style = $styleDonor.Style
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.013'
$ws.Range("D10").Style = $styleDonor.Style
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.48'
$ws.Range("D11").Style = $styleDonor.Style
$ws.Range("E11").Value = '  +4.33%  '
$ws.Range("D12").Value = '1.908.06'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.013'
$ws.Range("D13").Style = $styleDonor.Style
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.141'
$ws.Range("D14").Style = $styleDonor.Style
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.22'
$ws.Range("D15").Style = $styleDonor.Style
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06762'
$ws.Range("D17").Style = $styleDonor.Style
$ws.Range("E17").Value = '  +1.90%  '
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("D21").Value = '29.500.92'
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.624'
$ws.Range("D22").Style = $styleDonor.Style
$ws.Range("E22").Value = '  +1.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.77'
$ws.Range("D23").Style = $styleDonor.Style
$ws.Range("E23").Value = '  +2.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.183'
$ws.Range("D24").Style = $styleDonor.Style
$ws.Range("E24").Value = '  -0.71%  '
$ws.Range("D25").Value = '2.122.62'
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.374'
$ws.Range("D27").Style = $styleDonor.Style
$ws.Range("E27").Value = '  +5.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.03'
$ws.Range("D28").Style = $styleDonor.Style
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.109'
$ws.Range("D29").Style = $styleDonor.Style
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.87'
$ws.Range("D30").Style = $styleDonor.Style
$ws.Range("E30").Value = '  +2.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.024'
$ws.Range("D31").Style = $styleDonor.Style
$ws.Range("E31").Value = '  -4.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09533'
$ws.Range("D32").Style = $styleDonor.Style
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.521'
$ws.Range("D33").Style = $styleDonor.Style
$ws.Range("E33").Value = '  +2.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.568'
$ws.Range("D34").Style = $styleDonor.Style
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.389'
$ws.Range("D35").Style = $styleDonor.Style
$ws.Range("E35").Value = '  -2.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02269'
$ws.Range("D36").Style = $styleDonor.Style
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06096'
$ws.Range("D37").Style = $styleDonor.Style
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.175'
$ws.Range("D38").Style = $styleDonor.Style
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5974'
$ws.Range("D39").Style = $styleDonor.Style
$ws.Range("E39").Value = '  +1.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.988'
$ws.Range("D40").Style = $styleDonor.Style
$ws.Range("E40").Value = '  -3.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.76'
$ws.Range("D41").Style = $styleDonor.Style
$ws.Range("E41").Value = '  +6.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1855'
$ws.Range("D42").Style = $styleDonor.Style
$ws.Range("E42").Value = '  +1.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.284'
$ws.Range("D43").Style = $styleDonor.Style
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.407'
$ws.Range("D44").Style = $styleDonor.Style
$ws.Range("E44").Value = '  -4.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.56'
$ws.Range("D45").Style = $styleDonor.Style
$ws.Range("E45").Value = '  +4.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07624'
$ws.Range("D46").Style = $styleDonor.Style
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5576'
$ws.Range("D47").Style = $styleDonor.Style
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.940'
$ws.Range("D48").Style = $styleDonor.Style
$ws.Range("E48").Value = '  +0.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '115.84'
$ws.Range("D49").Style = $styleDonor.Style
$ws.Range("E49").Value = '  +2.61%  '

# Row 50 and 51: MXToken and Aave swap positions with updated values
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.56'
$ws.Range("D50").Style = $styleDonor.Style
$ws.Range("E50").Value = '  +1.87%  '

$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.416'
$ws.Range("D51").Style = $styleDonor.Style
$ws.Range("E51").Value = '  +2.91%  '
